$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = '0.489 (0.476 ± 0.010)'
$ws.Range("C2").Value = '00:04:59 (00:08:18 ± 00:02:34)'
$ws.Range("D2").Value = '00:00:01 (00:00:04 ± 00:00:02)'
$ws.Range("B3").Value = '0.516 (0.475 ± 0.018)'
$ws.Range("C3").Value = '00:01:05 (00:01:22 ± 00:00:11)'
$ws.Range("D3").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B4").Value = '0.487 (0.430 ± 0.034)'
$ws.Range("C4").Value = '00:00:43 (00:00:59 ± 00:00:16)'
$ws.Range("D4").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B5").Value = '0.474 (0.227 ± 0.148)'
$ws.Range("C5").Value = '00:05:05 (00:05:13 ± 00:00:03)'
$ws.Range("D5").Value = '00:00:01 (00:00:02 ± 00:00:02)'
$ws.Range("B6").Value = '0.533 (0.492 ± 0.026)'
$ws.Range("C6").Value = '00:04:56 (00:05:01 ± 00:00:02)'
$ws.Range("D6").Value = '00:00:00 (00:00:02 ± 00:00:01)'
$ws.Range("B7").Value = '0.469 (0.468 ± 0.001)'
$ws.Range("C7").Value = '00:05:03 (00:05:06 ± 00:00:02)'
$ws.Range("D7").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B9").Value = '0.496 (0.464 ± 0.022)'
$ws.Range("C9").Value = '00:05:00 (00:05:04 ± 00:00:03)'
$ws.Range("D9").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B10").Value = '0.456 (0.450 ± 0.005)'
$ws.Range("C10").Value = '00:04:29 (00:04:29 ± 00:00:00)'
$ws.Range("D10").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B11").Value = '0.370 (0.307 ± 0.045)'
$ws.Range("C11").Value = '00:05:06 (00:05:07 ± 00:00:00)'
$ws.Range("D11").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B12").Value = '0.506 (0.463 ± 0.028)'
$ws.Range("C12").Value = '00:02:03 (00:05:11 ± 00:02:05)'
$ws.Range("D12").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B13").Value = '0.146 (0.099 ± 0.033)'
$ws.Range("C13").Value = '00:00:08 (00:00:09 ± 00:00:00)'
$ws.Range("D13").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B14").Value = '0.477 (0.449 ± 0.019)'
$ws.Range("C14").Value = '00:01:44 (00:01:56 ± 00:00:07)'
$ws.Range("D14").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B15").Value = '0.495 (0.462 ± 0.016)'
$ws.Range("C15").Value = '00:01:22 (00:04:30 ± 00:00:55)'
$ws.Range("D15").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B16").Value = '0.515 (0.468 ± 0.025)'
$ws.Range("C16").Value = '00:05:05 (00:08:54 ± 00:06:57)'
$ws.Range("D16").Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range("B17").Value = '0.504 (0.460 ± 0.022)'
$ws.Range("C17").Value = '00:05:03 (00:05:51 ± 00:00:26)'
$ws.Range("D17").Value = '00:00:00 (00:00:00 ± 00:00:00)'